$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "TEst 2"
$ws.Range("C2").Value = "09:35"
$ws.Range("D2").Value = "Test"

# Row 3
$ws.Range("B3").Value = "Test"
$ws.Range("C3").Value = "09:35"
$ws.Range("D3").Value = "TEst 3"

# Row 4
$ws.Range("B4").Value = "TEst 3"
$ws.Range("C4").Value = "09:35"
$ws.Range("D4").Value = "TEst4"
$ws.Range("F4").Value = 1

# Row 5
$ws.Range("B5").Value = "TEst4"
$ws.Range("C5").Value = "09:35"
